# "add small animation start for adding souls" -- dev log update.
#
# Substantive changes:
#  1. TODO sheet: "add gold coins" and the related consumable-tablets task
#     are marked "canceled" (gold coins concept dropped).
#  2. TODO sheet: the "think to get rid of gold coins concept..." task is
#     marked "done", with a completion date.
#  3. Logs sheet: a new dated log entry is appended describing the long
#     gap and the work that resumed (adding the souls view/images).
#  4. The "Logs" tab becomes the active/selected sheet (was "TODO Before
#     0.0.1"), reflecting where work happened last.

$wb = $excel.ActiveWorkbook

$todo = $wb.Worksheets.Item("TODO Before 0.0.1")
$logs = $wb.Worksheets.Item("Logs")

# --- TODO Before 0.0.1 -----------------------------------------------

# "add gold coins" -> canceled
$todo.Range("C37").Value = "canceled"

# "gold coins could be used to buy the consumable items..." -> canceled
$todo.Range("C38").Value = "canceled"

# "think to get rid of gold coins concept - better sould to be used" -> done
$todo.Range("C48").Value = "done"
$todo.Range("D48").Value = 45673
$todo.Range("D48").NumberFormat = "m/d/yy"

# Leave the old selection on this sheet where the last edit happened.
$todo.Range("B42").Select()

# --- Logs --------------------------------------------------------------

$logs.Range("A70").Value = 45673
$logs.Range("A70").NumberFormat = "m/d/yy"
$logs.Range("B70").Value = "wow what an gap. Long hollidays. Decided to move forward whith what I have to finish project as demo game. Added the souls view and images."

# Make Logs the active sheet, with the cell below the new entry selected,
# matching where Excel leaves the cursor after typing a log line.
$logs.Range("B71").Select()
$logs.Activate()
